# Teste de busca data+sensação termica
# Shift the existing currency-quote log down by one "snapshot" (4 rows: Dollar,
# Euro, Iene, Yuan Chinês), insert a brand new snapshot at the top (rows 2-5)
# and re-append the oldest snapshot (which falls off the bottom) as new rows
# at the end of the sheet, bumping the running index in column A and growing
# the used range from A1:E105 to A1:E109.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 105
$numRows      = $lastDataRow - $firstDataRow + 1   # 104 existing data rows

# --- 1. Read the whole existing data block (A2:E105) into memory -----------
# Note: Range.Value2 returns a 1-based 2D array: $old[row,col]
$srcRange = $ws.Range("A$firstDataRow`:E$lastDataRow")
$old = $srcRange.Value2

# --- 2. Build the new data block (108 rows: 104 old + 4 new) ----------------
# Note: New-Object 'object[,]' returns a 0-based 2D array: $new[row,col]
$newNumRows = $numRows + 4
$new = New-Object 'object[,]' $newNumRows,5

# Brand-new snapshot inserted at the very top of the sheet
$topValues = @(
    @("Dollar",      "4,86",  "01:14", " quarta-feira, 14 de junho de 2023 "),
    @("Euro",         "5,24", "01:14", " quarta-feira, 14 de junho de 2023 "),
    @("Iene",         "0,035","01:14", " quarta-feira, 14 de junho de 2023 "),
    @("Yuan Chinês",  "0,68", "01:14", " quarta-feira, 14 de junho de 2023 ")
)

for ($i = 0; $i -lt 4; $i++) {
    $new[$i,1] = $topValues[$i][0]
    $new[$i,2] = $topValues[$i][1]
    $new[$i,3] = $topValues[$i][2]
    $new[$i,4] = $topValues[$i][3]
}

# Every other old row (1-based index $j within $old, j = 1..numRows) shifts
# down by 4 rows; this also naturally re-appends the previously-last 4 rows
# (the oldest snapshot) as the new final rows 106-109.
for ($j = 1; $j -le $numRows; $j++) {
    $new[$j+3,1] = $old[$j,2]
    $new[$j+3,2] = $old[$j,3]
    $new[$j+3,3] = $old[$j,4]
    $new[$j+3,4] = $old[$j,5]
}

# Column A is just the running index 0..(newNumRows-1)
for ($i = 0; $i -lt $newNumRows; $i++) {
    $new[$i,0] = $i
}

# Values such as "0,037" look like a number using a comma thousands
# separator (e.g. "123,456" -> 123456), so Excel's text/number inference
# would silently turn them into numbers. Prefix those with a quote so they
# are stored as genuine text, exactly like the source data. Remember which
# ones need the quote-prefix "undone" (style reset) after writing.
$forcedTextCells = @()
for ($i = 0; $i -lt $newNumRows; $i++) {
    for ($col = 1; $col -le 4; $col++) {
        $cellVal = $new[$i,$col]
        if ($cellVal -is [string] -and $cellVal -match '^-?\d+,\d{3}$') {
            $new[$i,$col] = "'" + $cellVal
            $forcedTextCells += , @($i, $col)
        }
    }
}

# --- 3. Write the new block back to the sheet --------------------------------
$lastNewRow = $firstDataRow + $newNumRows - 1   # 109
$dstRange = $ws.Range("A$firstDataRow`:E$lastNewRow")
$dstRange.Value2 = $new

# The quote-prefix trick marks the cell's style with quotePrefix="1" so it
# keeps showing as plain text; reset the style back to Normal (matching the
# original/unstyled cells in columns B-E) now that the value has "stuck".
foreach ($pair in $forcedTextCells) {
    $r = $firstDataRow + $pair[0]
    $c = $pair[1] + 1
    $ws.Cells.Item($r, $c).Style = "Normal"
}

# --- 4. Make sure the 4 brand-new rows (106-109) have the same column-A
#        style (bold, bordered, centered) as the rest of the index column ---
$styleSource = $ws.Cells.Item($lastDataRow, 1)   # A105, already styled
for ($r = $lastDataRow + 1; $r -le $lastNewRow; $r++) {
    $target = $ws.Cells.Item($r, 1)
    $val = $target.Value2
    $styleSource.Copy($target)
    $target.Value2 = $val
}

# --- 5. Refresh the sheet's reported dimension -------------------------------
$ws.Range("A1:E$lastNewRow").Select() | Out-Null
